$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column AB header "user_mobile" (style copied from B1 - header style)
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("AB1").Value = "user_mobile"

# ---------------------------------------------------------------------------
# 2. Column AB width
# ---------------------------------------------------------------------------
$ws.Columns("AB").ColumnWidth = 17.6

# ---------------------------------------------------------------------------
# 3. Data rows 2-30 : user_mobile phone numbers
#    style "7" (copied from E2) for most rows, style "1" (copied from B2) for
#    rows 3 and 4.
# ---------------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("AB2").PasteSpecial(-4122)
$ws.Range("AB2").Value = 9874563215

$ws.Range("B2").Copy()
$ws.Range("AB3").PasteSpecial(-4122)
$ws.Range("AB3").Value = 9874563216

$ws.Range("B2").Copy()
$ws.Range("AB4").PasteSpecial(-4122)
$ws.Range("AB4").Value = 9874563216

$ws.Range("E2").Copy()
$ws.Range("AB5").PasteSpecial(-4122)
$ws.Range("AB5").Value = 4525235325

$ws.Range("E2").Copy()
$ws.Range("AB6:AB30").PasteSpecial(-4122)
$ws.Range("AB6:AB30").Value = 9874563215

# ---------------------------------------------------------------------------
# 4. Rows 31-57 : empty cells, style "7"
# ---------------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("AB31:AB57").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Rows 58-63 : empty cells, style "2"
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("AB58:AB63").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Sheet view: drop the frozen/scrolled topLeftCell and move the selection
#    from C32 to F3.
# ---------------------------------------------------------------------------
$null = $ws.Range("F3").Select()
